# Modify how models are picked
# Fill in actual results (Actual/Correct columns) for the games that had
# already been played, and append the next slate of predicted games.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fill in Actual (F) and Correct (G) columns for rows 222-231 ---
# These games had already been predicted; now we know what actually happened.
$results = @(
    @{ Row = 222; Actual = "Brandon Wheat Kings";    Correct = 1 },
    @{ Row = 223; Actual = "Prince Albert Raiders";  Correct = 0 },
    @{ Row = 224; Actual = "Saskatoon Blades";        Correct = 0 },
    @{ Row = 225; Actual = "Lethbridge Hurricanes";  Correct = 1 },
    @{ Row = 226; Actual = "Red Deer Rebels";         Correct = 0 },
    @{ Row = 227; Actual = "Kamloops Blazers";        Correct = 0 },
    @{ Row = 228; Actual = "Everett Silvertips";      Correct = 1 },
    @{ Row = 229; Actual = "Portland Winterhawks";    Correct = 1 },
    @{ Row = 230; Actual = "Tri-City Americans";      Correct = 0 },
    @{ Row = 231; Actual = "Victoria Royals";         Correct = 1 }
)

foreach ($r in $results) {
    $ws.Range("F" + $r.Row).Value = $r.Actual
    $ws.Range("G" + $r.Row).Value = $r.Correct
}

# --- Append new games (Sat, Feb 15, 2025) as rows 232-241 ---
$newGames = @(
    @{ Row = 232; GameID = 1021778; Home = "Brandon Wheat Kings";   Away = "Moose Jaw Warriors";    Prediction = "Brandon Wheat Kings" },
    @{ Row = 233; GameID = 1021782; Home = "Prince Albert Raiders"; Away = "Calgary Hitmen";         Prediction = "Calgary Hitmen" },
    @{ Row = 234; GameID = 1021786; Home = "Swift Current Broncos"; Away = "Saskatoon Blades";       Prediction = "Saskatoon Blades" },
    @{ Row = 235; GameID = 1021779; Home = "Edmonton Oil Kings";    Away = "Lethbridge Hurricanes";  Prediction = "Lethbridge Hurricanes" },
    @{ Row = 236; GameID = 1021781; Home = "Medicine Hat Tigers";   Away = "Regina Pats";            Prediction = "Medicine Hat Tigers" },
    @{ Row = 237; GameID = 1021783; Home = "Prince George Cougars"; Away = "Everett Silvertips";     Prediction = "Everett Silvertips" },
    @{ Row = 238; GameID = 1021780; Home = "Kelowna Rockets";       Away = "Kamloops Blazers";       Prediction = "Kelowna Rockets" },
    @{ Row = 239; GameID = 1021784; Home = "Seattle Thunderbirds";  Away = "Portland Winterhawks";   Prediction = "Seattle Thunderbirds" },
    @{ Row = 240; GameID = 1021785; Home = "Spokane Chiefs";        Away = "Tri-City Americans";     Prediction = "Tri-City Americans" },
    @{ Row = 241; GameID = 1021787; Home = "Vancouver Giants";      Away = "Victoria Royals";        Prediction = "Victoria Royals" }
)

foreach ($g in $newGames) {
    $ws.Range("A" + $g.Row).Value = $g.GameID
    $ws.Range("B" + $g.Row).Value = "Sat, Feb 15, 2025"
    $ws.Range("C" + $g.Row).Value = $g.Home
    $ws.Range("D" + $g.Row).Value = $g.Away
    $ws.Range("E" + $g.Row).Value = $g.Prediction
}

# --- Update the view so the newly entered rows are in focus ---
[void]$ws.Range("I232").Select()
try {
    $excel.ActiveWindow.ScrollRow = 211
} catch {
}
